$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 98.0
$ws.Range("B3").Value = 47.0
$ws.Range("B4").Value = 73.0
